$p = $ppt.ActivePresentation

# Slide 10
$s = $p.Slides.Item(10)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Protect NATW statutory basis (CRS 23-52-105) from DEI misclassification"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Prioritize retention as most cost-effective enrollment strategy"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Frame Indigenous education through statutory/sovereign obligations, not DEI language"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Invest in AI Institute and place-based experiential learning"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Qualify online expansion: Indigenous niche only (NATW moat), not generic degrees"

# Slide 3
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: High (5/5) | Trend: Negative"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Trump administration (2025–2029) reducing federal HE funding; 120 TRIO programs terminated"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "DEI programs under HIGH scrutiny — executive order targeting DEI in accreditation (Apr 2025)"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Tribal education funding VOLATILE: 109% increase Sept 2025, but FY2026 proposes 24% cuts"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Colorado FY 2025–26: `$38.4M increase (far less than `$95M requested); 3.5% tuition cap"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: Reframe Indigenous programs through statutory obligations (CRS 23-52-105) and cultural preservation (legally safe)"

# Slide 4
$s = $p.Slides.Item(4)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: High (5/5) | Trend: Negative"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Colorado shifts costs to students via tuition rather than state appropriations"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Rising tuition sensitivity; students increasingly price-conscious and comparison-shopping"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Durango housing crisis — major hidden barrier for student attendance AND faculty recruitment"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Native American tuition waiver revenue impact (~37% of students at zero tuition)"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: Healthcare/nursing programs (strong regional employer demand)"

# Slide 5
$s = $p.Slides.Item(5)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: Medium-High (4/5) | Trend: Mixed"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Declining college-going rates nationally and in Colorado"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Career outcome expectations dominant ('What job will I get?')"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Indigenous education opportunity IS REAL (166 tribes, 37% waiver, underserved nationally)"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "First-generation students (43%) need targeted support systems"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: Indigenous education leadership — reframe through statutory obligations (CRS 23-52-105), not DEI"

# Slide 6
$s = $p.Slides.Item(6)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: High (4/5) | Trend: Rapidly Changing"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "AI disruption transforming pedagogy, assessment, and student expectations"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Online graduate market SATURATED — ASU, SNHU, Western Governors dominate (`$50M+ marketing)"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "FLC has NO online brand nationally; ~25 online courses (~10% of offerings)"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Passive video lectures becoming obsolete; AI-enabled adaptive learning replacing them"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: AI Institute partnerships and curriculum integration"

# Slide 7
$s = $p.Slides.Item(7)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: High (4/5) | Trend: Deteriorating"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Title VI scrutiny — 50+ universities under investigation for race-conscious programs"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Native American Tuition Waiver has DISTINCT legal basis (CRS 23-52-105, since 1911)"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "HLC accreditation: federal pressure on DEI standards, but HLC offers flexibility"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Trump administration revising Title IX regulations (definitions, due process in flux)"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: NATW defensible under Title VI (statutory basis per CRS 23-52-105, not voluntary DEI)"

# Slide 8
$s = $p.Slides.Item(8)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Impact: Medium (3/5) | Trend: Negative"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "Southwest Colorado wildfire risk increasing — smoke impacts air quality and outdoor activities"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Colorado River basin under long-term drought stress; water rights contentious"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Snowpack variability affects regional economy (ski, rafting, outdoor recreation)"
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "Outdoor recreation brand is FLC strength but CLIMATE-VULNERABLE"
$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Opportunity: Proactive sustainability initiatives to build brand beyond compliance"

# Slide 9
$s = $p.Slides.Item(9)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = ""
$tr.Paragraphs(1, 1).Text = "Highest impact: Political (5/5), Economic (5/5)"
$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "High impact: Social (4/5), Technological (4/5), Legal (4/5)"
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "Key risks: Federal DEI policy, tribal waiver vulnerability, Durango housing crisis"
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "Key opportunity: Indigenous education (statutorily grounded), AI Institute, experiential learning"
